$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing B/C values for rows 2-5 (use .Formula because .Value
# is unreliable for reads in this runtime).
$b2 = $ws.Cells.Item(2, 2).Formula
$c2 = $ws.Cells.Item(2, 3).Formula
$b3 = $ws.Cells.Item(3, 2).Formula
$c3 = $ws.Cells.Item(3, 3).Formula
$b4 = $ws.Cells.Item(4, 2).Formula
$c4 = $ws.Cells.Item(4, 3).Formula
$b5 = $ws.Cells.Item(5, 2).Formula
$c5 = $ws.Cells.Item(5, 3).Formula

# Row 2 becomes what used to be row 5
$ws.Cells.Item(2, 2).Formula = $b5
$ws.Cells.Item(2, 3).Formula = $c5
$ws.Cells.Item(2, 4).Value = ""

# Row 3 becomes what used to be row 4
$ws.Cells.Item(3, 2).Formula = $b4
$ws.Cells.Item(3, 3).Formula = $c4
$ws.Cells.Item(3, 4).Value = ""

# Row 4 becomes what used to be row 3
$ws.Cells.Item(4, 2).Formula = $b3
$ws.Cells.Item(4, 3).Formula = $c3
$ws.Cells.Item(4, 4).Value = ""

# Row 5 becomes what used to be row 2
$ws.Cells.Item(5, 2).Formula = $b2
$ws.Cells.Item(5, 3).Formula = $c2
$ws.Cells.Item(5, 4).Value = ""
